# FTUX: Content - Added new columns to individually control when Passive
# Events and Happy Hour are allowed.
#
# The "gameSettings" table on the global_settings sheet gains two new
# trailing columns: [enablePassiveEventsAtRun] and [enableHappyHourAtRun]3.
# Both get header text in row 4 and a value of 4 in row 5, matching the
# formatting of the preceding [enableShareButtonsAtRun] column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item(1)

# Grow the table by two columns (B4:N5 -> B4:P5).
$lo.Resize($ws.Range("B4:P5"))

# New header names (row 4) - setting the header cell text renames the
# corresponding table column.
$ws.Range("O4").Value = "[enablePassiveEventsAtRun]"
$ws.Range("P4").Value = "[enableHappyHourAtRun]3"

# New data values (row 5).
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 4

# Match the look of the last existing column (N) for both the header and
# data rows.
$ws.Range("N4").Copy()
$ws.Range("O4:P4").PasteSpecial(-4122)

$ws.Range("N5").Copy()
$ws.Range("O5:P5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Leave the selection where the author ended up (just past the new table).
$ws.Range("Q5").Select()
